# Slide 16, "Content Placeholder 2" shape: collapse the
# ": 1h per week, 28h hours in " run sequence into a single
# ": 2*1*14h = 28h in " explanation run.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(16)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$oldFragment = ": 1h per week, 28h hours in "
$newFragment = ": 2*1*14h = 28h in "

$full = $tr.Text
$idx  = $full.IndexOf($oldFragment)

$target = $tr.Characters($idx + 1, $oldFragment.Length)
$target.Text = $newFragment
